$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.280.42'
$ws.Range('E2').Value = '  +5.17%  '
$ws.Range('D3').Value = '3.464.43'
$ws.Range('E3').Value = '  +5.22%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '186.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '547.35'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.614'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.48%  '
$ws.Range('D8').Value = '3.455.73'
$ws.Range('E8').Value = '  +5.15%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.638'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.10'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.147'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +10.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000274'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.47'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.39%  '
$ws.Range('D15').Value = '4.026.57'
$ws.Range('E15').Value = '  +5.85%  '
$ws.Range('D16').Value = '3.465.84'
$ws.Range('E16').Value = '  +5.39%  '
$ws.Range('D17').Value = '67.648.09'
$ws.Range('E17').Value = '  +6.02%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.120'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.86%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.29'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.68%  '
$ws.Range('E21').Value = '  +5.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '406.28'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.93'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.08%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.74'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.19%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.89'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.20'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +9.52%  '
$ws.Range('E28').Value = '  +2.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.75'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.65'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.20'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '681.52'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.88'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.93%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.66'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.110'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.12'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = '0.0₃0828'
$ws.Range('E37').Value = '  +17.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '38.69'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.403'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.78%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  +14.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.37'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +21.69%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.133'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('D45').Value = '3.051.87'
$ws.Range('E45').Value = '  +4.42%  '
$ws.Range('E46').Value = '  +10.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0419'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.29'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.73'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.11%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +14.89%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +9.54%  '
